# This script reproduces the data/format update described by the diff:
#  - column B (prediction score) values for rows 2-19 are overwritten with the
#    real classifier scores (they were all placeholder 1's before)
#  - the "Row"/label column (A) together with the header row (A1:C1) gets its
#    text format re-applied (the source pipeline re-writes this style on every
#    run, which is why the style table grows between successive outputs)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated prediction scores for column B, rows 2-19 ---
$newValues = @{
    2  = -0.64365691826914784
    3  = -0.73306600376966546
    4  =  0.25567047870162707
    5  = -0.49024892940747389
    6  = -0.45276723073732228
    7  =  0.30027566638045844
    8  =  0.4176294025220848
    9  = -0.13279277492446351
    10 =  0.26882819376633194
    11 =  0.24272154879670271
    12 =  0.36146099292739109
    13 =  0.40173077737204554
    14 =  0.0082978031163261079
    15 = -0.29784373368653139
    16 = -0.30090737181108551
    17 = -0.28870227443576901
    18 =  0.44003202845268063
    19 = -0.13548426709168027
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

# --- Re-apply the text ("@") number format to the row/label column, as the
#     originating tool does on every write, bumping the style table. ---
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A19").NumberFormat = "@"
